# Update the numeric weight/bias tables on the Layer0 and Layer1 sheets
# (Network4.xlsx) with newly-generated values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Layer0")
$ws.Range("B2").Value = -4.072171871222555
$ws.Range("C2").Value = -8.670785174811794
$ws.Range("D2").Value = -1.45889724731158
$ws.Range("E2").Value = 3.592368819678708
$ws.Range("B3").Value = -1.437424209417388
$ws.Range("C3").Value = -2.729848430930387
$ws.Range("D3").Value = -1.439438419614557
$ws.Range("E3").Value = 1.746985084067611
$ws.Range("B4").Value = -4.330335562621634
$ws.Range("C4").Value = 8.172728619266813
$ws.Range("D4").Value = 1.244843975562981
$ws.Range("E4").Value = -2.10006424343622
$ws.Range("B5").Value = 2.232486114250346
$ws.Range("C5").Value = -0.8992039525888637
$ws.Range("D5").Value = -0.7477647011036654
$ws.Range("E5").Value = -3.233102461719949
$ws.Range("B6").Value = 2.967362062903506
$ws.Range("C6").Value = -5.788525128429651
$ws.Range("D6").Value = 3.062927996034302
$ws.Range("E6").Value = -2.143104742091968
$ws.Range("B7").Value = 0.1399808187949514
$ws.Range("C7").Value = -0.5869834479577366
$ws.Range("D7").Value = 0.9302772084627584
$ws.Range("E7").Value = 7.206543998779834

$ws = $wb.Worksheets.Item("Layer1")
$ws.Range("B2").Value = -15.3251005650841
$ws.Range("C2").Value = 0.2147204893769506
$ws.Range("D2").Value = -11.01661000811683
$ws.Range("E2").Value = -1.039309771046578
$ws.Range("F2").Value = 0.4745177555307734
$ws.Range("B3").Value = 5.101917402969105
$ws.Range("C3").Value = 9.335444508167916
$ws.Range("D3").Value = -6.142790574187227
$ws.Range("E3").Value = -6.584404340445493
$ws.Range("F3").Value = -9.28423520033923
$ws.Range("B4").Value = -4.719634936671806
$ws.Range("C4").Value = -6.684705795183908
$ws.Range("D4").Value = 10.01923504110058
$ws.Range("E4").Value = -9.374111467761026
$ws.Range("F4").Value = 1.709578136882074
$ws.Range("B5").Value = 10.15944491586395
$ws.Range("C5").Value = -9.916586944108111
$ws.Range("D5").Value = -8.078926594891261
$ws.Range("E5").Value = -7.828683213980202
$ws.Range("F5").Value = 7.916993007091226
$ws.Range("B6").Value = 8.450296127014388
$ws.Range("C6").Value = -8.679110283866637
$ws.Range("D6").Value = 5.414968949804011
$ws.Range("E6").Value = 6.601693718153819
$ws.Range("F6").Value = -9.090553166238742
